# parametric bootstrap from Trosvik et al 2012 poles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add dec/inc values (E3, F3), matching the plain numeric style
#     already used by the neighboring G3/H3 cells ---
$ws.Range("G3").Copy()
$ws.Range("E3:F3").PasteSpecial(-4122)
$ws.Range("E3").Value = 351.3
$ws.Range("F3").Value = 39.9

# --- kept_by_auth column (V): give the new cells the same plain style
#     already used across the row (same as T7, e.g.) before filling them in ---
$ws.Range("T7").Copy()
$ws.Range("V7:V24").PasteSpecial(-4122)

$ws.Range("V7").Value = "n"
$ws.Range("V8").Value = "n"
$ws.Range("V9").Value = "n"
$ws.Range("V10").Value = "n"
$ws.Range("V11").Value = "n"
$ws.Range("V12").Value = "n"
$ws.Range("V13").Value = "n"
$ws.Range("V14").Value = "n"
$ws.Range("V15").Value = "n"

$ws.Range("V16").Value = "y"
$ws.Range("V17").Value = "y"
$ws.Range("V18").Value = "y"
$ws.Range("V19").Value = "y"
$ws.Range("V20").Value = "y"
$ws.Range("V21").Value = "y"
$ws.Range("V22").Value = "y"
$ws.Range("V23").Value = "y"
$ws.Range("V24").Value = "y"

# --- Highlight cells (cyan fill) for bootstrap-affected values ---
$cyan = 16776960
$ws.Range("H12").Interior.Color = $cyan
$ws.Range("H12").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("H17").PasteSpecial(-4122)
$ws.Range("G18").PasteSpecial(-4122)
$ws.Range("H18").PasteSpecial(-4122)
$ws.Range("G19").PasteSpecial(-4122)
$ws.Range("H19").PasteSpecial(-4122)
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("G21").PasteSpecial(-4122)
$ws.Range("H21").PasteSpecial(-4122)

# --- Column V width widened to fit "kept_by_auth" values ---
$ws.Range("V1").ColumnWidth = 11.5

"done"
